$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record for Espinaca / Femacal de La Calera was inserted
# at row 177, pushing the existing records (rows 177..283) down by one row
# (to rows 178..284). Copy whole rows (columns A:R) from the bottom up so
# that no row is overwritten before it has been read.
for ($r = 284; $r -ge 178; $r--) {
    $srcRow = $r - 1
    $src = $ws.Range("A" + $srcRow + ":R" + $srcRow).Value2
    $ws.Range("A" + $r + ":R" + $r).Value = $src
}

# Row 284 is brand new (the sheet used to end at 283), so its "Fecha" cell
# (D284) needs the same date format the rest of column D already carries.
$ws.Range("D284").NumberFormat = $ws.Range("D283").NumberFormat

# Write the new record into row 177 (same Mercado/Region/Categoria/Calidad/
# Unidad/Origen/Clasificacion as its neighbours, new Fecha + price figures).
$ws.Range("D177").Value = 44606
$ws.Range("J177").Value = 80
$ws.Range("K177").Value = 4000
$ws.Range("L177").Value = 4500
$ws.Range("M177").Value = 4281
$ws.Range("P177").Value = 1427
